$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "'23.76"
$ws.Range("D4").Value = "'5.247"
$ws.Range("D5").Value = "'0.05815"
$ws.Range("D6").Value = "'6.466"
$ws.Range("D7").Value = "'3.227"
$ws.Range("D8").Value = "'0.8083"
$ws.Range("D9").Value = "'0.8871"
$ws.Range("D10").Value = "'0.1392"
$ws.Range("D11").Value = "'0.07094"
$ws.Range("D12").Value = "'0.03106"
$ws.Range("D13").Value = "'0.03045"
$ws.Range("D14").Value = "'0.09326"
$ws.Range("D15").Value = "'3.843"
$ws.Range("D16").Value = "'0.001536"
$ws.Range("D17").Value = "'0.04703"
$ws.Range("D18").Value = "'0.0006013"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006152"
$ws.Range("D20").Value = "'0.001257"
$ws.Range("D22").Value = "'0.00008706"
$ws.Range("D24").Value = "'2.158"
$ws.Range("D25").Value = "'0.3180"
$ws.Range("D28").Value = "'0.0002329"
$ws.Range("D40").Value = "'0.03788"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1054"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002504"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.006274"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("D44").Value = "'0.007833"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.5352"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "'0.003249"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
